# Increment the "Förändrad" date (column C) by one day for all data rows (2-130)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 130
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
